# Apply updated cryptocurrency price/volume data to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.337.04"
$ws.Range("E2").Value = "  +0.62%  "
$ws.Range("D3").Value = "1.667.62"
$ws.Range("E3").Value = "  +0.84%  "
$ws.Range("D4").Value = "'1.010"
$ws.Range("E4").Value = "  +0.43%  "
$ws.Range("D5").Value = "'220.69"
$ws.Range("E5").Value = "  +1.35%  "
$ws.Range("D6").Value = "'0.5316"
$ws.Range("E6").Value = "  +0.04%  "
$ws.Range("D7").Value = "'1.010"
$ws.Range("E7").Value = "  +0.41%  "
$ws.Range("D8").Value = "'0.2653"
$ws.Range("E8").Value = "  +1.20%  "
$ws.Range("D9").Value = "'0.06363"
$ws.Range("E9").Value = "  +0.52%  "
$ws.Range("D10").Value = "'20.83"
$ws.Range("E10").Value = "  +2.17%  "
$ws.Range("D11").Value = "'0.07854"
$ws.Range("E11").Value = "  +0.67%  "
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "1.687.43"
$ws.Range("E12").Value = "  +2.61%  "
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").Value = "'4.509"
$ws.Range("E13").Value = "  -0.22%  "
$ws.Range("D14").Value = "1.896.28"
$ws.Range("E14").Value = "  +0.81%  "
$ws.Range("D15").Value = "'0.5587"
$ws.Range("E15").Value = "  +1.78%  "
$ws.Range("D16").Value = "0.0₅8150"
$ws.Range("E16").Value = "  -0.31%  "
$ws.Range("D17").Value = "'65.88"
$ws.Range("E17").Value = "  +0.74%  "
$ws.Range("D18").Value = "26.343.55"
$ws.Range("E18").Value = "  +0.70%  "
$ws.Range("D19").Value = "'1.010"
$ws.Range("E19").Value = "  +0.45%  "
$ws.Range("D20").Value = "'4.707"
$ws.Range("E20").Value = "  +2.43%  "
$ws.Range("D21").Value = "'197.05"
$ws.Range("E21").Value = "  +3.10%  "
$ws.Range("D22").Value = "'10.27"
$ws.Range("E22").Value = "  +1.79%  "
$ws.Range("D23").Value = "'6.040"
$ws.Range("E23").Value = "  +0.55%  "
$ws.Range("D24").Value = "'1.011"
$ws.Range("E24").Value = "  +0.37%  "
$ws.Range("D25").Value = "'145.39"
$ws.Range("E25").Value = "  +0.09%  "
$ws.Range("D26").Value = "'0.1220"
$ws.Range("E26").Value = "  -0.30%  "
$ws.Range("D27").Value = "'7.232"
$ws.Range("E27").Value = "  +0.52%  "
$ws.Range("D28").Value = "'16.20"
$ws.Range("E28").Value = "  +1.36%  "
$ws.Range("D29").Value = "'1.509"
$ws.Range("E29").Value = "  +2.52%  "
$ws.Range("D30").Value = "'0.05888"
$ws.Range("E30").Value = "  +2.60%  "
$ws.Range("D31").Value = "'1.284"
$ws.Range("E31").Value = "  +0.84%  "
$ws.Range("D32").Value = "'3.540"
$ws.Range("E32").Value = "  -0.11%  "
$ws.Range("E33").Value = "  +1.86%  "
$ws.Range("D34").Value = "'1.604"
$ws.Range("E34").Value = "  +0.68%  "
$ws.Range("D35").Value = "'0.9632"
$ws.Range("E35").Value = "  +1.63%  "
$ws.Range("D36").Value = "'2.832"
$ws.Range("E36").Value = "  +0.90%  "
$ws.Range("D37").Value = "'2.438"
$ws.Range("E37").Value = "  +0.68%  "
$ws.Range("D38").Value = "'0.5803"
$ws.Range("E38").Value = "  +1.20%  "
$ws.Range("E39").Value = "  +0.80%  "
$ws.Range("B40").Value = "FraxShare"
$ws.Range("C40").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D40").Value = "'5.938"
$ws.Range("E40").Value = "  +2.70%  "
$ws.Range("B41").Value = "Maker"
$ws.Range("C41").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D41").Value = "1.074.09"
$ws.Range("E41").Value = "  +3.49%  "
$ws.Range("D42").Value = "'0.8637"
$ws.Range("E42").Value = "  +1.75%  "
$ws.Range("D43").Value = "'1.010"
$ws.Range("E43").Value = "  +0.45%  "
$ws.Range("D44").Value = "'102.89"
$ws.Range("E44").Value = "  -1.00%  "
$ws.Range("D45").Value = "1.807.11"
$ws.Range("E45").Value = "  +0.73%  "
$ws.Range("D46").Value = "'58.24"
$ws.Range("E46").Value = "  +2.58%  "
$ws.Range("D47").Value = "0.0₈107"
$ws.Range("E47").Value = "  +1.36%  "
$ws.Range("E48").Value = "  +0.93%  "
$ws.Range("D49").Value = "'0.4413"
$ws.Range("E49").Value = "  +1.31%  "
$ws.Range("D50").Value = "'7.980"
$ws.Range("E50").Value = "  +1.16%  "
$ws.Range("D51").Value = "'0.05152"
$ws.Range("E51").Value = "  -0.07%  "
